{"js": "// Rename table header labels:\n//   \"\u0418\u043c\u044f \u0440\u0435\u0433\u0438\u0441\u0442\u0440\u0430\" -> \"\u0420\u0435\u0433\u0438\u0441\u0442\u0440\"\n//   \"\u0418\u043c\u044f \u043f\u043e\u043b\u044f\"     -> \"\u041f\u043e\u043b\u0435\"\n//   \"Reset\"        -> \"\u0417\u043d\u0430\u0447\u0435\u043d\u0438\u0435\"\nconst replacements = [\n  [\"\u0418\u043c\u044f \u0440\u0435\u0433\u0438\u0441\u0442\u0440\u0430\", \"\u0420\u0435\u0433\u0438\u0441\u0442\u0440\"],\n  [\"\u0418\u043c\u044f \u043f\u043e\u043b\u044f\", \"\u041f\u043e\u043b\u0435\"],\n  [\"Reset\", \"\u0417\u043d\u0430\u0447\u0435\u043d\u0438\u0435\"],\n];\n\nfor (const [search, replacement] of replacements) {\n  const results = context.document.body.search(search, {\n    matchCase: true,\n    matchWholeWord: true,\n  });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Rename table header labels:\n#   \"\u0418\u043c\u044f \u0440\u0435\u0433\u0438\u0441\u0442\u0440\u0430\" -> \"\u0420\u0435\u0433\u0438\u0441\u0442\u0440\"\n#   \"\u0418\u043c\u044f \u043f\u043e\u043b\u044f\"     -> \"\u041f\u043e\u043b\u0435\"\n#   \"Reset\"        -> \"\u0417\u043d\u0430\u0447\u0435\u043d\u0438\u0435\"\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"\u0418\u043c\u044f \u0440\u0435\u0433\u0438\u0441\u0442\u0440\u0430\"; Replace = \"\u0420\u0435\u0433\u0438\u0441\u0442\u0440\" },\n    @{ Find = \"\u0418\u043c\u044f \u043f\u043e\u043b\u044f\"; Replace = \"\u041f\u043e\u043b\u0435\" },\n    @{ Find = \"Reset\"; Replace = \"\u0417\u043d\u0430\u0447\u0435\u043d\u0438\u0435\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    # wdFindContinue=1, Replace:=wdReplaceAll=2, MatchWholeWord=$true\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
